# "Generate Report for Handoff"
#
# The handoff-report generator just finished producing a new handoff
# package for file 27b28b76-0401-45e1-a793-72c236bdad2c, so its
# "Latest Handoff Datetime" column (H) gets stamped with the timestamp
# of the newly generated xlf, on both the zh-cn and the de-de status
# sheets (row 4 in each table is the 27b28b76 file).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("H4").Value = "2016-10-27 07:37:51"
$wsDeDe.Range("H4").Value = "2016-10-27 07:38:02"
